$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table is being extended with a new "2022" column (M), which mirrors
# the formatting of the existing "2021" column (L) and repeats its values,
# except for the year header itself.

# Copy column L's formatting (borders/fonts/number formats) into the new
# column M first, then overwrite the values with the 2022 figures.
$ws.Range("L3:L11").Copy()
$ws.Range("M3:M11").PasteSpecial(-4122)

$ws.Range("M4").Value = 2022
$ws.Range("M5").Value = 0.86
$ws.Range("M6").Value = 1.07
$ws.Range("M7").Value = 25.27
$ws.Range("M8").Value = 14
$ws.Range("M9").Value = 0.12
$ws.Range("M10").Value = 21.74
$ws.Range("M11").Value = 9.4600000000000009

# The active selection in the workbook moved to N6.
$ws.Range("N6").Select()
